$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" conversion summary text (cell A1) ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$rangeA1 = $wsHoja1.Range("A1")
$oldText = $rangeA1.Value()
$newText = $oldText.Replace("1000 Bs = 7.1 = 28360.54 pesos", "1000 Bs = 7.1 = 28423.3 pesos")
$newText = $newText.Replace("28360.54 pesos = 7.03 = 959.88 Bs", "28423.3 pesos = 7.08 = 977.9 Bs")
$rangeA1.Value = $newText

# --- Update the "tasas" rate table (cells N10, O10, N12, O12) ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 140.8
$wsTasas.Range("O10").Value = 4002
$wsTasas.Range("N12").Value = 4013.95
$wsTasas.Range("O12").Value = 138.099
